$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 4.425865000000001
$ws.Range("H2").Value = 13.277595
$ws.Range("I2").Value = 0.4619841037548157
$ws.Range("J2").Value = 0.4696223785602887
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 44.13164066666667
$ws.Range("N2").Value = 132.394922
$ws.Range("O2").Value = 0.4415399811720331
$ws.Range("P2").Value = 0.4562856844211927
$ws.Range("Q2").Value = 195.3206838191767
$ws.Range("R2").Value = 1757.88615437259
$ws.Range("S2").Value = 0.2039844524736799
$ws.Range("T2").Value = 0.2142819684208898

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 4.425865000000001
$ws.Range("H3").Value = 13.277595
$ws.Range("I3").Value = 0.4619841037548157
$ws.Range("J3").Value = 0.4696223785602887
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 14.93259333333333
$ws.Range("N3").Value = 44.79778
$ws.Range("O3").Value = 0.1494015830739255
$ws.Range("P3").Value = 0.1543910098595022
$ws.Range("Q3").Value = 66.08964219323336
$ws.Range("R3").Value = 594.8067797391001
$ws.Range("S3").Value = 0.06902115645595812
$ws.Range("T3").Value = 0.07250547327854438

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 4.425865000000001
$ws.Range("H4").Value = 13.277595
$ws.Range("I4").Value = 0.4619841037548157
$ws.Range("J4").Value = 0.4696223785602887
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.61024133333333
$ws.Range("N4").Value = 40.830724
$ws.Range("O4").Value = 0.1361713639304118
$ws.Range("P4").Value = 0.1407189532975654
$ws.Range("Q4").Value = 60.23709075875334
$ws.Range("R4").Value = 542.13381682878
$ws.Range("S4").Value = 0.06290900552246212
$ws.Range("T4").Value = 0.06608476955611685

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.425865000000001
$ws.Range("H5").Value = 13.277595
$ws.Range("I5").Value = 0.4619841037548157
$ws.Range("J5").Value = 0.4696223785602887
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.584752
$ws.Range("N5").Value = 52.754256
$ws.Range("O5").Value = 0.1759366057935712
$ws.Range("P5").Value = 0.1818121982434553
$ws.Range("Q5").Value = 77.82773841048001
$ws.Range("R5").Value = 700.4496456943201
$ws.Range("S5").Value = 0.08127991514520733
$ws.Range("T5").Value = 0.08538307699036624

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.425865000000001
$ws.Range("H6").Value = 13.277595
$ws.Range("I6").Value = 0.4619841037548157
$ws.Range("J6").Value = 0.4696223785602887
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 9.690137500000001
$ws.Range("N6").Value = 19.380275
$ws.Range("O6").Value = 0.09695046603005844
$ws.Range("P6").Value = 0.06679215417828435
$ws.Range("Q6").Value = 42.88724040643751
$ws.Range("R6").Value = 257.323442438625
$ws.Range("S6").Value = 0.04478957415750826
$ws.Range("T6").Value = 0.03136709031437142

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.686805000000001
$ws.Range("H7").Value = 14.060415
$ws.Range("I7").Value = 0.4892217470254038
$ws.Range("J7").Value = 0.4973103589802793
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 44.13164066666667
$ws.Range("N7").Value = 132.394922
$ws.Range("O7").Value = 0.4415399811720331
$ws.Range("P7").Value = 0.4562856844211927
$ws.Range("Q7").Value = 206.8363941347367
$ws.Range("R7").Value = 1861.52754721263
$ws.Range("S7").Value = 0.216010960970546
$ws.Range("T7").Value = 0.2269155975170658

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.686805000000001
$ws.Range("H8").Value = 14.060415
$ws.Range("I8").Value = 0.4892217470254038
$ws.Range("J8").Value = 0.4973103589802793
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 14.93259333333333
$ws.Range("N8").Value = 44.79778
$ws.Range("O8").Value = 0.1494015830739255
$ws.Range("P8").Value = 0.1543910098595022
$ws.Range("Q8").Value = 69.98615309763335
$ws.Range("R8").Value = 629.8753778787001
$ws.Range("S8").Value = 0.07309050347978682
$ws.Range("T8").Value = 0.07678024853655685

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.686805000000001
$ws.Range("H9").Value = 14.060415
$ws.Range("I9").Value = 0.4892217470254038
$ws.Range("J9").Value = 0.4973103589802793
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.61024133333333
$ws.Range("N9").Value = 40.830724
$ws.Range("O9").Value = 0.1361713639304118
$ws.Range("P9").Value = 0.1407189532975654
$ws.Range("Q9").Value = 63.78854713227334
$ws.Range("R9").Value = 574.09692419046
$ws.Range("S9").Value = 0.06661799255686811
$ws.Range("T9").Value = 0.06998099317974141

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.686805000000001
$ws.Range("H10").Value = 14.060415
$ws.Range("I10").Value = 0.4892217470254038
$ws.Range("J10").Value = 0.4973103589802793
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 17.584752
$ws.Range("N10").Value = 52.754256
$ws.Range("O10").Value = 0.1759366057935712
$ws.Range("P10").Value = 0.1818121982434553
$ws.Range("Q10").Value = 82.41630359736
$ws.Range("R10").Value = 741.7467323762401
$ws.Range("S10").Value = 0.08607201365205071
$ws.Range("T10").Value = 0.09041708957544647

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.686805000000001
$ws.Range("H11").Value = 14.060415
$ws.Range("I11").Value = 0.4892217470254038
$ws.Range("J11").Value = 0.4973103589802793
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 9.690137500000001
$ws.Range("N11").Value = 19.380275
$ws.Range("O11").Value = 0.09695046603005844
$ws.Range("P11").Value = 0.06679215417828435
$ws.Range("Q11").Value = 45.41578488568751
$ws.Range("R11").Value = 272.494709314125
$ws.Range("S11").Value = 0.04743027636615225
$ws.Range("T11").Value = 0.03321643017146875

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.467454
$ws.Range("H12").Value = 0.9349080000000001
$ws.Range("I12").Value = 0.04879414921978045
$ws.Range("J12").Value = 0.03306726245943202
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 44.13164066666667
$ws.Range("N12").Value = 132.394922
$ws.Range("O12").Value = 0.4415399811720331
$ws.Range("P12").Value = 0.4562856844211927
$ws.Range("Q12").Value = 20.629511956196
$ws.Range("R12").Value = 123.777071737176
$ws.Range("S12").Value = 0.02154456772780724
$ws.Range("T12").Value = 0.01508811848323715

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.467454
$ws.Range("H13").Value = 0.9349080000000001
$ws.Range("I13").Value = 0.04879414921978045
$ws.Range("J13").Value = 0.03306726245943202
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 14.93259333333333
$ws.Range("N13").Value = 44.79778
$ws.Range("O13").Value = 0.1494015830739255
$ws.Range("P13").Value = 0.1543910098595022
$ws.Range("Q13").Value = 6.980300484040002
$ws.Range("R13").Value = 41.88180290424
$ws.Range("S13").Value = 0.007289923138180545
$ws.Range("T13").Value = 0.005105288044400915

$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 0.467454
$ws.Range("H14").Value = 0.9349080000000001
$ws.Range("I14").Value = 0.04879414921978045
$ws.Range("J14").Value = 0.03306726245943202
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 13.61024133333333
$ws.Range("N14").Value = 40.830724
$ws.Range("O14").Value = 0.1361713639304118
$ws.Range("P14").Value = 0.1407189532975654
$ws.Range("Q14").Value = 6.362161752232001
$ws.Range("R14").Value = 38.172970513392
$ws.Range("S14").Value = 0.006644365851081541
$ws.Range("T14").Value = 0.004653190561707153

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 0.467454
$ws.Range("H15").Value = 0.9349080000000001
$ws.Range("I15").Value = 0.04879414921978045
$ws.Range("J15").Value = 0.03306726245943202
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 17.584752
$ws.Range("N15").Value = 52.754256
$ws.Range("O15").Value = 0.1759366057935712
$ws.Range("P15").Value = 0.1818121982434553
$ws.Range("Q15").Value = 8.220062661407999
$ws.Range("R15").Value = 49.320375968448
$ws.Range("S15").Value = 0.008584676996313205
$ws.Range("T15").Value = 0.006012031677642624

$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 0.467454
$ws.Range("H16").Value = 0.9349080000000001
$ws.Range("I16").Value = 0.04879414921978045
$ws.Range("J16").Value = 0.03306726245943202
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 9.690137500000001
$ws.Range("N16").Value = 19.380275
$ws.Range("O16").Value = 0.09695046603005844
$ws.Range("P16").Value = 0.06679215417828435
$ws.Range("Q16").Value = 4.529693534925
$ws.Range("R16").Value = 18.1187741397
$ws.Range("S16").Value = 0.004730615506397928
$ws.Range("T16").Value = 0.002208633692444178
